$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data for rows 2-5 (columns A:T), matching the target OOXML state.
$data = @(
    @("FAPs","Wnt5a","Fzd6","ECs",   3,1,4.619088000000001,13.857264,1,1,2,1,25.081504,50.163008,0.670801886161467,0.5799583902495139,115.853674148352,695.1220448901121,0.670801886161467,0.5799583902495139),
    @("FAPs","Wnt5a","Fzd6","FAPs",  3,1,4.619088000000001,13.857264,1,1,2,0.6666666666666666,0.3140343333333334,0.942103,0.008398811455623189,0.01089210079525609,1.450552220688,13.054969986192,0.008398811455623189,0.01089210079525609),
    @("FAPs","Wnt5a","Fzd6","Neutro",3,1,4.619088000000001,13.857264,1,1,3,1,11.39945533333333,34.198366,0.3048770974345635,0.3953835721837834,52.65508733673601,473.895786030624,0.3048770974345635,0.3953835721837834),
    @("FAPs","Wnt5a","Fzd6","sCs",   3,1,4.619088000000001,13.857264,1,1,2,1,0.5953364999999999,1.190673,0.01592220494834624,0.01376593677144639,2.749911683112,16.499470098672,0.01592220494834624,0.01376593677144639)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($c = 1; $c -le $vals.Length; $c++) {
        $ws.Cells.Item($row, $c).Value = $vals[$c - 1]
    }
}
